$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (current row 2 "Gobernatura 18" data shifts down to row 3, etc.)
$ws.Rows.Item(2).Insert()

# Fill in the new row 2 with the new "Gobernatura 17" entry
$ws.Range("A2").Value = "Gobernatura 17"
$ws.Range("B2").Value = "gb_17"
$ws.Range("C2").Value = "#dda15e"

# Update selection to mirror the authored workbook state
$ws.Range("C2").Select()
